$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title "Attendance Taking System": bump font size to 18pt (sz/szCs=36)
#    on both the paragraph mark run-props (pPr/rPr) and the run itself.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Font.Size = 18
$titlePara.Range.Font.SizeBi = 18

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: remove it from the end of the "Once all
#    checks..." paragraph and add it to the empty paragraph right before it.
# ---------------------------------------------------------------------------

# Find the empty paragraph immediately preceding the "Once all checks..."
# paragraph, and the "Once all checks..." paragraph itself.
$goBackHost = $null
$onceParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Once all checks have been done*") {
        $onceParaIndex = $i
        break
    }
}
$goBackHostIndex = $onceParaIndex - 1

# Remove the old "_GoBack" bookmark (hidden from Bookmarks.Count, but still
# addressable by name) before we touch the paragraph's runs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Add the new "_GoBack" bookmark at the empty paragraph before it.
$hostPara = $d.Paragraphs($goBackHostIndex)
$d.Bookmarks.Add("_GoBack", $hostPara.Range)

# ---------------------------------------------------------------------------
# 3) Merge the runs of the "Once all checks..." paragraph (drops the
#    spellStart/spellEnd proofErr wrapper around "its").
# ---------------------------------------------------------------------------
$oncePara = $d.Paragraphs($onceParaIndex)
$r = $d.Range($oncePara.Range.Start, $oncePara.Range.End - 1)
# Force a real mutation even though the visible text is unchanged, since a
# same-text assignment is treated as a no-op by the engine.
$r.Text = "placeholder"
$oncePara2 = $d.Paragraphs($onceParaIndex)
$r2 = $d.Range($oncePara2.Range.Start, $oncePara2.Range.End - 1)
$apos = [char]0x2019
$r2.Text = "Once all checks have been done and its attendance taken, it will lastly send a notification to the User" + $apos + "s phone stating its attendance have been taken."

# ---------------------------------------------------------------------------
# 4) "StudentDetails" -> "Student Details" (drops spellStart/spellEnd).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(StudentDetails)") | Out-Null
$rng.Text = "(Student Details)"

# ---------------------------------------------------------------------------
# 5) "Javascript" -> "JavaScript" (drops spellStart/spellEnd).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("using Javascript (") | Out-Null
$rng.Text = "using JavaScript ("

# ---------------------------------------------------------------------------
# 6) "A comments page" -> "A comment page" (gramStart/gramEnd moves, but the
#    visible content change is comments -> comment).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("A comments page") | Out-Null
$rng.Text = "A comment page"

Write-Output "done"
